# clean_EMX_func for removal of unallowed chars
# This edit adds the missing "ID" attribute row for the rd_bb_contribution
# entity on the "attributes" sheet, and the corresponding "ID" column on the
# "rd_bb_contribution" sheet, matching the pattern already used by every
# other entity (rd_reg_accessibility, rd_scientific publications, etc).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "attributes" sheet: insert a new "ID" attribute row right before the
#    first existing rd_bb_contribution attribute row (currently row 120).
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("attributes")

$ws.Rows.Item(120).Insert()

$ws.Cells.Item(120, 1).Value = "ID"
$ws.Cells.Item(120, 2).Value = "ID"
$ws.Cells.Item(120, 3).Value = " "
$ws.Cells.Item(120, 4).Value = "rd_bb_contribution"
$ws.Cells.Item(120, 8).Value = "'true"
$ws.Cells.Item(120, 8).Style = "Normal"
$ws.Cells.Item(120, 10).Value = "'true"
$ws.Cells.Item(120, 10).Style = "Normal"

# ---------------------------------------------------------------------
# 2) "rd_bb_contribution" sheet: insert a new "ID" column at the front so
#    it matches every other entity sheet (ID always comes first).
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("rd_bb_contribution")

$ws2.Columns.Item(1).Insert()
$ws2.Cells.Item(1, 1).Value = "ID"
